# Auto-generated Excel COM-interop script
# Applies numeric value updates (and a few cell additions/removals)
# to the Halicarnassus_Profits workbook sheets, per the target diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 8455.833000000001
$ws.Range("I2").Value = 16711.666
$ws.Range("K2").Value = 16711.666
$ws.Range("M2").Value = -16598.666
$ws.Range("H76").Value = 300
$ws.Range("I76").Value = 300
$ws.Range("K76").Value = 300
$ws.Range("M76").Value = 15
$ws.Range("H79").Value = 300
$ws.Range("I79").Value = 300
$ws.Range("K79").Value = 300
$ws.Range("M79").Value = 792
$ws.Range("H97").Value = 200
$ws.Range("J97").Value = 200
$ws.Range("L97").Value = 600
$ws.Range("N97").Value = -1592
$ws.Range("H103").Value = 3825.4
$ws.Range("I103").Value = 10001
$ws.Range("J103").Value = 2281.5
$ws.Range("K103").Value = 30003
$ws.Range("L103").Value = 6844.5
$ws.Range("M103").Value = -29417
$ws.Range("N103").Value = -8016.5
$ws.Range("H104").Value = 661.1667
$ws.Range("I104").Value = 661.1667
$ws.Range("K104").Value = 1983.5001
$ws.Range("M104").Value = -236.5001
$ws.Range("H118").Value = 2000
$ws.Range("I118").Value = 2000
$ws.Range("K118").Value = 6000
$ws.Range("M118").Value = -4343
$ws.Range("H138").Value = 4160
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4160
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12480
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -22760
$ws.Range("H141").Value = 3000
$ws.Range("I141").Value = 3250
$ws.Range("K141").Value = 9750
$ws.Range("M141").Value = -4570

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 20000950
$ws.Range("H26").Value = 2023.75
$ws.Range("I26").Value = 2023.75
$ws.Range("K26").Value = 2023.75
$ws.Range("M26").Value = -1693.75
$ws.Range("H122").Value = 1957.5
$ws.Range("I122").Value = 1994.1666
$ws.Range("K122").Value = 5982.4998
$ws.Range("M122").Value = -3532.4998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3371.8215
$ws.Range("I86").Value = 2023.3889
$ws.Range("J86").Value = 5799
$ws.Range("K86").Value = 2023.3889
$ws.Range("L86").Value = 5799
$ws.Range("M86").Value = -900.3888999999999
$ws.Range("N86").Value = -8045
$ws.Range("H89").Value = 3371.8215
$ws.Range("I89").Value = 2023.3889
$ws.Range("J89").Value = 5799
$ws.Range("K89").Value = 10116.9445
$ws.Range("L89").Value = 28995
$ws.Range("M89").Value = -4500.9445
$ws.Range("N89").Value = -40227
$ws.Range("H95").Value = 7208
$ws.Range("J95").Value = 7208
$ws.Range("L95").Value = 7208
$ws.Range("N95").Value = -12700

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 108.5
$ws.Range("I19").Value = 129.16667
$ws.Range("J19").Value = 46.5
$ws.Range("K19").Value = 129.16667
$ws.Range("L19").Value = 46.5
$ws.Range("M19").Value = 40.83332999999999
$ws.Range("N19").Value = -386.5
$ws.Range("H24").Value = 108.5
$ws.Range("I24").Value = 129.16667
$ws.Range("J24").Value = 46.5
$ws.Range("K24").Value = 129.16667
$ws.Range("L24").Value = 46.5
$ws.Range("M24").Value = 40.83332999999999
$ws.Range("N24").Value = -386.5
$ws.Range("H31").Value = 8574.5
$ws.Range("I31").Value = 3140
$ws.Range("J31").Value = 9828.615
$ws.Range("K31").Value = 3140
$ws.Range("L31").Value = 9828.615
$ws.Range("M31").Value = -2845
$ws.Range("N31").Value = -10418.615
$ws.Range("H34").Value = 8574.5
$ws.Range("I34").Value = 3140
$ws.Range("J34").Value = 9828.615
$ws.Range("K34").Value = 3140
$ws.Range("L34").Value = 9828.615
$ws.Range("M34").Value = -2938
$ws.Range("N34").Value = -10232.615
$ws.Range("H62").Value = 5600.077
$ws.Range("I62").Value = 4749.8335
$ws.Range("J62").Value = 6328.857
$ws.Range("K62").Value = 4749.8335
$ws.Range("L62").Value = 6328.857
$ws.Range("M62").Value = -4125.8335
$ws.Range("N62").Value = -7576.857
$ws.Range("H65").Value = 5600.077
$ws.Range("I65").Value = 4749.8335
$ws.Range("J65").Value = 6328.857
$ws.Range("K65").Value = 23749.1675
$ws.Range("L65").Value = 31644.285
$ws.Range("M65").Value = -20629.1675
$ws.Range("N65").Value = -37884.285
$ws.Range("H88").Value = 17500
$ws.Range("J88").Value = 17500
$ws.Range("L88").Value = 17500
$ws.Range("N88").Value = -18312
$ws.Range("H91").Value = 17500
$ws.Range("J91").Value = 17500
$ws.Range("L91").Value = 17500
$ws.Range("N91").Value = -20308
$ws.Range("H99").Value = 1879.3334
$ws.Range("J99").Value = 2025
$ws.Range("L99").Value = 2025
$ws.Range("N99").Value = -5021
$ws.Range("H126").Value = 1879.3334
$ws.Range("J126").Value = 2025
$ws.Range("L126").Value = 6075
$ws.Range("N126").Value = -11015

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 412.125
$ws.Range("I38").Value = 415.6
$ws.Range("J38").Value = 394.75
$ws.Range("K38").Value = 1246.8
$ws.Range("L38").Value = 1184.25
$ws.Range("M38").Value = -899.8000000000002
$ws.Range("N38").Value = -1878.25
$ws.Range("H109").Value = 6844.1763
$ws.Range("J109").Value = 13697.5
$ws.Range("L109").Value = 41092.5
$ws.Range("N109").Value = -43172.5
$ws.Range("H113").Value = 602.8570999999999
$ws.Range("I113").Value = 254
$ws.Range("J113").Value = 1475
$ws.Range("K113").Value = 762
$ws.Range("L113").Value = 4425
$ws.Range("M113").Value = 1408
$ws.Range("N113").Value = -8765
$ws.Range("H114").Value = 745.1667
$ws.Range("I114").Value = 294.2
$ws.Range("K114").Value = 882.5999999999999
$ws.Range("M114").Value = 2371.4
$ws.Range("H137").Value = 5362.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 5362.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 16087.5
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -26287.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6141.143
$ws.Range("I122").Value = 5747
$ws.Range("J122").Value = 6666.6665
$ws.Range("K122").Value = 17241
$ws.Range("L122").Value = 19999.9995
$ws.Range("M122").Value = -14791
$ws.Range("N122").Value = -24899.9995
$ws.Range("H141").Value = 67500
$ws.Range("J141").Value = 67500
$ws.Range("L141").Value = 67500
$ws.Range("N141").Value = -77860

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 950
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 950
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -655
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 950
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -843
$ws.Range("N27").ClearContents()
$ws.Range("H93").Value = 1703.6666
$ws.Range("I93").Value = 2409.6667
$ws.Range("K93").Value = 2409.6667
$ws.Range("M93").Value = -1161.6667
$ws.Range("H100").Value = 7984.9
$ws.Range("I100").Value = 1450
$ws.Range("K100").Value = 1450
$ws.Range("M100").Value = -909
$ws.Range("H132").Value = 3150
$ws.Range("I132").Value = 2687.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8062.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5532.5
$ws.Range("N132").Value = -20060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6331
$ws.Range("J9").Value = 10666
$ws.Range("L9").Value = 10666
$ws.Range("N9").Value = -10946
$ws.Range("H110").Value = 25000
$ws.Range("J110").Value = 25000
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180
$ws.Range("H135").Value = 39905
$ws.Range("J135").Value = 39905
$ws.Range("L135").Value = 39905
$ws.Range("N135").Value = -50045
$ws.Range("H140").Value = 45000
$ws.Range("J140").Value = 45000
$ws.Range("L140").Value = 45000
$ws.Range("N140").Value = -55360
$ws.Range("H141").Value = 299997.5
$ws.Range("J141").Value = 99995
$ws.Range("L141").Value = 99995
$ws.Range("N141").Value = -110355
